$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - PACOTE PRÉ-OPERATÓRIO PEDIÁTRICO OTORRINO
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 12
$ws.Range("F2").Value = 6
$ws.Range("G2").Value = 3
$ws.Range("I2").Value = 5
$ws.Range("J2").Value = 1
$ws.Range("M2").Value = 1

# Row 3 - PACOTE PRÉ-OPERATÓRIO PEDIÁTRICO CIRURGIA GERAL
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 10
$ws.Range("G3").Value = 6
$ws.Range("H3").Value = 3
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 7
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3

# Row 5 - ADENOIDECTOMIA PEDIÁTRICO
$ws.Range("C5").Value = 2
$ws.Range("G5").Value = 2
$ws.Range("I5").Value = 1

# Row 6 - AMIGDALECTOMIA- PEDIATRICO (label fixed, removed space before hyphen)
$ws.Range("A6").Value = "AMIGDALECTOMIA- PEDIATRICO"
$ws.Range("C6").Value = 1
$ws.Range("I6").Value = 1

# Row 7 - AMIGDALECTOMIA COM ADENOIDECTOMIA - PEDIATRICO
$ws.Range("C7").Value = 5
$ws.Range("F7").Value = 6
$ws.Range("G7").Value = 1
$ws.Range("I7").Value = 3
$ws.Range("M7").Value = 1

# Row 10 - HERNIOPLASTIA INGUINAL (BILATERAL) - PEDIATRICO
$ws.Range("E10").Value = 0
$ws.Range("G10").Value = 1
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 1
$ws.Range("M10").Value = 0

# Row 11 - HERNIOPLASTIA UMBILICAL - PEDIATRICO
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0

# Row 12 - ORQUIDOPEXIA BILATERAL - PEDIATRICO
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1
$ws.Range("K12").Value = 0
$ws.Range("M12").Value = 1

# Row 14 - CORRECAO DE HIPOSPADIA (1º TEMPO) - PEDIATRICO
$ws.Range("K14").Value = 0

# Row 16 - POSTECTOMIA - PEDIATRICO
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 4
$ws.Range("H16").Value = 2
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 1

# Row 17 - TOTAL
$ws.Range("B17").Value = 2
$ws.Range("C17").Value = 20
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 29
$ws.Range("G17").Value = 18
$ws.Range("H17").Value = 5
$ws.Range("I17").Value = 10
$ws.Range("J17").Value = 1
$ws.Range("K17").Value = 13
$ws.Range("L17").Value = 2
$ws.Range("M17").Value = 7
